$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper behaviour note: this COM-interop host coalesces adjacent runs that
# end up with identical run properties whenever text is replaced in place
# (e.g. via Find.Execute or Range.Text=). However, if a sub-range is given a
# run-level formatting change (even one that is immediately reverted back to
# match its neighbours) the run boundary it creates is preserved. We use a
# tiny Bold toggle (set then immediately unset) on the remainder range
# starting at each desired split point to carve the paragraph's single run
# into the multiple runs required, without altering the final formatting.
# ---------------------------------------------------------------------------

function Split-RunAt([int]$absPos, [int]$rangeEnd) {
    $r = $d.Range($absPos, $rangeEnd)
    $r.Bold = 1
    $r.Bold = 0
}

# ---------------------------------------------------------------------------
# Edit 1: "... показывает высокие результаты. Использование ..."
#      -> "... показал высокие результаты быстродействия. Использование ..."
# ---------------------------------------------------------------------------
$found1 = $d.Content.Find.Execute("показывает высокие результаты.", $true, $false, $false, $false, $false, $true, 1, $false, "показал высокие результаты быстродействия.", 2)

$full = $d.Content.Text
$marker1 = "Далее были проведены испытания быстродействия эмуляторов. Разработанный эмулятор показал высокие результаты быстродействия. Использование набирающего популярность средства проектирования ПО "
$sentStart1 = $full.IndexOf($marker1)
$sentEnd1 = $sentStart1 + $marker1.Length

# internal split offsets (relative to $sentStart1): after "показ", after "ал", after "высокие результаты", after " быстродействия"
foreach ($off in @(122, 107, 88, 86)) {
    Split-RunAt ($sentStart1 + $off) $sentEnd1
}

# ---------------------------------------------------------------------------
# Edit 2: ". Разработанный эмулятор показывает низкие результаты"
#      -> ". Разработанный эмулятор показал низкие результаты"
# ---------------------------------------------------------------------------
$found2 = $d.Content.Find.Execute("показывает низкие результаты", $true, $false, $false, $false, $false, $true, 1, $false, "показал низкие результаты", 2)

$full2 = $d.Content.Text
$marker2 = ". Разработанный эмулятор показал низкие результаты"
$sentStart2 = $full2.IndexOf($marker2)
$sentEnd2 = $sentStart2 + $marker2.Length

# internal split offsets: after "показ", after "ал"
foreach ($off in @(32, 30)) {
    Split-RunAt ($sentStart2 + $off) $sentEnd2
}

# ---------------------------------------------------------------------------
# Edit 3: add <w:lang w:val="en-US"/> to the empty paragraph right after the
# "... эмулируемой системы." paragraph (paraId 02B21DCD).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "") {
        $prevText = ""
        if ($i -gt 1) { $prevText = $d.Paragraphs($i - 1).Range.Text }
        if ($prevText.Contains("эмулируемой системы")) {
            $p.Range.LanguageID = "en-US"
        }
    }
}

# ---------------------------------------------------------------------------
# Edit 4: ", по результатам которой выбрано средство для реализации эмулятора ядра "
#      -> ", по результатам которого выбрано средство для реализации эмулятора ядра "
# ---------------------------------------------------------------------------
$found4 = $d.Content.Find.Execute("по результатам которой выбрано", $true, $false, $false, $false, $false, $true, 1, $false, "по результатам которого выбрано", 2)

$full4 = $d.Content.Text
$marker4 = ", по результатам которого выбрано средство для реализации эмулятора ядра "
$sentStart4 = $full4.IndexOf($marker4)
$sentEnd4 = $sentStart4 + $marker4.Length

# internal split offset: after "которо"
foreach ($off in @(23)) {
    Split-RunAt ($sentStart4 + $off) $sentEnd4
}

Write-Host "found1=$found1 found2=$found2 found4=$found4"
